$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Task₃" header to "Task₃ Parcial 1"
$ws.Range("E1").Value = "Task₃ Parcial 1"

# Fill in "Task₃ Parcial 1" grades for each student (rows 2-32).
# Most students get a full numeric grade of 5; a few get a partial
# (decimal) grade that must be stored as text, matching the source data.
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 5
$ws.Range("E6").Value = 5

# E7 (BULA VERGARA MARIA PAULA) -> "2.3" stored as text
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.3"
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").Value = 5
$ws.Range("E9").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("E11").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("E13").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("E16").Value = 5
$ws.Range("E17").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("E19").Value = 5
$ws.Range("E20").Value = 5
$ws.Range("E21").Value = 5

# E22 (PEREZ THERAN ERLINDA) -> "3.4" stored as text
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.4"
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E23").Value = 5

# D24 / E24 (QUINONEZ CERVANTES MARIANA) -> Task2 grade 5, Task3 "3.3" as text
$ws.Range("D24").Value = 5
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.3"
$ws.Range("E24").NumberFormat = "General"

$ws.Range("E25").Value = 5
$ws.Range("E26").Value = 5
$ws.Range("E27").Value = 5
$ws.Range("E28").Value = 5
$ws.Range("E29").Value = 5
$ws.Range("E30").Value = 5
$ws.Range("E31").Value = 5
$ws.Range("E32").Value = 5

# Update the active cell selection to match the author's final cursor position
$null = $ws.Range("E26").Select()
